$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 324.25
$ws.Range("I11").Value = 324.25
$ws.Range("K11").Value = 324.25
$ws.Range("M11").Value = -184.25
$ws.Range("H40").Value = 900
$ws.Range("I40").Value = 900
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 900
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -725
$ws.Range("N40").ClearContents()
$ws.Range("H113").Value = 3187.2222
$ws.Range("I113").Value = 2583.6
$ws.Range("J113").Value = 3941.75
$ws.Range("K113").Value = 2583.6
$ws.Range("L113").Value = 3941.75
$ws.Range("M113").Value = 670.4000000000001
$ws.Range("N113").Value = -10449.75
$ws.Range("H135").Value = 698.4545000000001
$ws.Range("I135").Value = 698.4545000000001
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6286.0905
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3751.0905
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 2516.8462
$ws.Range("J137").Value = 3411
$ws.Range("L137").Value = 10233
$ws.Range("N137").Value = -15333
$ws.Range("H138").Value = 1209.25
$ws.Range("I138").Value = 418.5
$ws.Range("K138").Value = 1255.5
$ws.Range("M138").Value = 3884.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 3600
$ws.Range("J10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("N10").Value = -5340
$ws.Range("H32").Value = 9397.177
$ws.Range("I32").Value = 8109.5
$ws.Range("K32").Value = 8109.5
$ws.Range("M32").Value = -7822.5
$ws.Range("H45").Value = 1849.4286
$ws.Range("I45").Value = 1845.1666
$ws.Range("J45").Value = 1875
$ws.Range("K45").Value = 1845.1666
$ws.Range("L45").Value = 1875
$ws.Range("M45").Value = -1468.1666
$ws.Range("N45").Value = -2629
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240
$ws.Range("H86").Value = 314314
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 314314
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 314314
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -316686
$ws.Range("H89").Value = 314314
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 314314
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 942942
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -954798
$ws.Range("H92").Value = 20250
$ws.Range("I92").Value = 10500
$ws.Range("K92").Value = 10500
$ws.Range("M92").Value = -8004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3684.8
$ws.Range("I20").Value = 2106
$ws.Range("K20").Value = 2106
$ws.Range("M20").Value = -1859
$ws.Range("H51").Value = 40000
$ws.Range("I51").Value = 40000
$ws.Range("K51").Value = 40000
$ws.Range("M51").Value = -39509
$ws.Range("H134").Value = 7079.448
$ws.Range("I134").Value = 7517.5
$ws.Range("J134").Value = 4976.8
$ws.Range("K134").Value = 22552.5
$ws.Range("L134").Value = 14930.4
$ws.Range("M134").Value = -20017.5
$ws.Range("N134").Value = -20000.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 29633.223
$ws.Range("J95").Value = 29633.223
$ws.Range("L95").Value = 29633.223
$ws.Range("N95").Value = -35125.223

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1643.4
$ws.Range("J5").Value = 490
$ws.Range("L5").Value = 1470
$ws.Range("N5").Value = -1694
$ws.Range("H17").Value = 576.2308
$ws.Range("J17").Value = 798.8889
$ws.Range("L17").Value = 2396.6667
$ws.Range("N17").Value = -2734.6667
$ws.Range("H68").Value = 2003
$ws.Range("J68").Value = 2003
$ws.Range("L68").Value = 6009
$ws.Range("N68").Value = -7631
$ws.Range("H71").Value = 2003
$ws.Range("J71").Value = 2003
$ws.Range("L71").Value = 18027
$ws.Range("N71").Value = -26139
$ws.Range("H107").Value = 568.25
$ws.Range("J107").Value = 1375
$ws.Range("L107").Value = 4125
$ws.Range("N107").Value = -7965
$ws.Range("H113").Value = 677.8
$ws.Range("J113").Value = 799.75
$ws.Range("L113").Value = 2399.25
$ws.Range("N113").Value = -6739.25
$ws.Range("H131").Value = 1468
$ws.Range("I131").Value = 937.25
$ws.Range("K131").Value = 2811.75
$ws.Range("M131").Value = 2228.25
$ws.Range("H132").Value = 1144.25
$ws.Range("I132").Value = 800
$ws.Range("J132").Value = 1488.5
$ws.Range("K132").Value = 7200
$ws.Range("L132").Value = 13396.5
$ws.Range("M132").Value = -4670
$ws.Range("N132").Value = -18456.5
$ws.Range("H135").Value = 1643.4
$ws.Range("J135").Value = 490
$ws.Range("L135").Value = 4410
$ws.Range("N135").Value = -9480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 389.05884
$ws.Range("I2").Value = 481
$ws.Range("K2").Value = 481
$ws.Range("M2").Value = -368
$ws.Range("H132").Value = 1855.5454
$ws.Range("I132").Value = 1855.5454
$ws.Range("K132").Value = 5566.6362
$ws.Range("M132").Value = -3036.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 39997
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 39997
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 39997
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -42243
$ws.Range("H83").Value = 39997
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 39997
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 119991
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -131223
$ws.Range("H99").Value = 5000
$ws.Range("I99").Value = 5000
$ws.Range("K99").Value = 5000
$ws.Range("M99").Value = -2005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 20005
$ws.Range("J14").Value = 20005
$ws.Range("L14").Value = 20005
$ws.Range("N14").Value = -20341
